# Update Fonds de solidarite - Volet 1 data with 2022-06-21 values.
# For each affected row, column C (nombre_aides) and column E (montant_total)
# are updated to their new totals. Column D (nombre_entreprises) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 6;   C = 20797;  E = 360751379  },
    @{ Row = 8;   C = 1053;   E = 91488760   },
    @{ Row = 24;  C = 35712;  E = 132422179  },
    @{ Row = 38;  C = 7252;   E = 58569613   },
    @{ Row = 58;  C = 394;    E = 34766328   },
    @{ Row = 70;  C = 15736;  E = 24686069   },
    @{ Row = 92;  C = 409331; E = 1597674047 },
    @{ Row = 93;  C = 209679; E = 1310315886 },
    @{ Row = 94;  C = 94245;  E = 919346057  },
    @{ Row = 96;  C = 17329;  E = 797818686  },
    @{ Row = 104; C = 135305; E = 272661062  },
    @{ Row = 114; C = 3805;   E = 9123175    },
    @{ Row = 120; C = 55;     E = 2597532    },
    @{ Row = 141; C = 80477;  E = 280751076  },
    @{ Row = 144; C = 24421;  E = 202173459  },
    @{ Row = 176; C = 28903;  E = 263366691  },
    @{ Row = 179; C = 635;    E = 54145745   }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
